$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "Grunnskolelærere"
$ws.Range("C23").Value = 52497
$ws.Range("D23").Value = 52497
$ws.Range("C24").Value = 52773
$ws.Range("D24").Value = 52459
$ws.Range("E24").Value = 314
$ws.Range("C25").Value = 53233
$ws.Range("D25").Value = 52487
$ws.Range("E25").Value = 746
$ws.Range("C26").Value = 53749
$ws.Range("D26").Value = 52346
$ws.Range("E26").Value = 1403
$ws.Range("C27").Value = 54321
$ws.Range("D27").Value = 52084
$ws.Range("E27").Value = 2237
$ws.Range("C28").Value = 54931
$ws.Range("D28").Value = 51638
$ws.Range("E28").Value = 3294
$ws.Range("C29").Value = 55547
$ws.Range("D29").Value = 51095
$ws.Range("E29").Value = 4453
$ws.Range("C30").Value = 56193
$ws.Range("D30").Value = 50601
$ws.Range("E30").Value = 5593
$ws.Range("C31").Value = 56819
$ws.Range("D31").Value = 50166
$ws.Range("E31").Value = 6653
$ws.Range("C32").Value = 57406
$ws.Range("D32").Value = 49727
$ws.Range("E32").Value = 7679
$ws.Range("C33").Value = 57928
$ws.Range("D33").Value = 49404
$ws.Range("E33").Value = 8524
$ws.Range("C34").Value = 58391
$ws.Range("D34").Value = 49125
$ws.Range("E34").Value = 9266
$ws.Range("C35").Value = 58778
$ws.Range("D35").Value = 48904
$ws.Range("E35").Value = 9874
$ws.Range("C36").Value = 59083
$ws.Range("D36").Value = 48741
$ws.Range("E36").Value = 10342
$ws.Range("C37").Value = 59331
$ws.Range("D37").Value = 48739
$ws.Range("E37").Value = 10592
$ws.Range("C38").Value = 59522
$ws.Range("D38").Value = 48844
$ws.Range("E38").Value = 10678
$ws.Range("C39").Value = 59689
$ws.Range("D39").Value = 49013
$ws.Range("E39").Value = 10676
$ws.Range("C40").Value = 59849
$ws.Range("D40").Value = 49230
$ws.Range("E40").Value = 10618
$ws.Range("C41").Value = 60004
$ws.Range("D41").Value = 49495
$ws.Range("E41").Value = 10509
$ws.Range("C42").Value = 60182
$ws.Range("D42").Value = 49814
$ws.Range("E42").Value = 10369
$ws.Range("C43").Value = 60374
$ws.Range("D43").Value = 50150
$ws.Range("E43").Value = 10224
$ws.Range("A44").Value = "Faglærere"
$ws.Range("C44").Value = 16600
$ws.Range("D44").Value = 16600
$ws.Range("C45").Value = 17578
$ws.Range("D45").Value = 16622
$ws.Range("E45").Value = 956
$ws.Range("C46").Value = 18647
$ws.Range("D46").Value = 16681
$ws.Range("E46").Value = 1966
$ws.Range("C47").Value = 19746
$ws.Range("D47").Value = 16724
$ws.Range("E47").Value = 3022
$ws.Range("C48").Value = 20847
$ws.Range("D48").Value = 16749
$ws.Range("E48").Value = 4098
$ws.Range("C49").Value = 21949
$ws.Range("D49").Value = 16780
$ws.Range("E49").Value = 5169
$ws.Range("C50").Value = 23031
$ws.Range("D50").Value = 16801
$ws.Range("E50").Value = 6230
$ws.Range("C51").Value = 24103
$ws.Range("D51").Value = 16800
$ws.Range("E51").Value = 7303
$ws.Range("C52").Value = 25145
$ws.Range("D52").Value = 16746
$ws.Range("E52").Value = 8400
$ws.Range("C53").Value = 26173
$ws.Range("D53").Value = 16669
$ws.Range("E53").Value = 9503
$ws.Range("C54").Value = 27167
$ws.Range("D54").Value = 16591
$ws.Range("E54").Value = 10576
$ws.Range("C55").Value = 28127
$ws.Range("D55").Value = 16535
$ws.Range("E55").Value = 11592
$ws.Range("C56").Value = 29061
$ws.Range("D56").Value = 16495
$ws.Range("E56").Value = 12567
$ws.Range("C57").Value = 29976
$ws.Range("D57").Value = 16480
$ws.Range("E57").Value = 13496
$ws.Range("C58").Value = 30886
$ws.Range("D58").Value = 16440
$ws.Range("E58").Value = 14447
$ws.Range("C59").Value = 31782
$ws.Range("D59").Value = 16386
$ws.Range("E59").Value = 15396
$ws.Range("C60").Value = 32658
$ws.Range("D60").Value = 16345
$ws.Range("E60").Value = 16313
$ws.Range("C61").Value = 33531
$ws.Range("D61").Value = 16350
$ws.Range("E61").Value = 17181
$ws.Range("C62").Value = 34404
$ws.Range("D62").Value = 16381
$ws.Range("E62").Value = 18024
$ws.Range("C63").Value = 35275
$ws.Range("D63").Value = 16417
$ws.Range("E63").Value = 18858
$ws.Range("C64").Value = 36130
$ws.Range("D64").Value = 16463
$ws.Range("E64").Value = 19667
